# Workbook / worksheet references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update the date in A1 (45308 -> 45309, i.e. one day later) ---
$ws.Range("A1").Value = 45309

# --- Step 2: update the price in D44 ---
$ws.Range("D44").Value = 43783.243

# --- Refresh the merged-cell regions (re-merge in the order that reproduces
#     the new mergeCells ordering seen after the edit) ---
$ws.Range("A10:D10").UnMerge()
$ws.Range("A10:D10").Merge()

$ws.Range("A11:D11").UnMerge()
$ws.Range("A11:D11").Merge()

$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:D1").Merge()

$ws.Range("B43:C43").UnMerge()
$ws.Range("B43:C43").Merge()

$ws.Range("B44:C44").UnMerge()
$ws.Range("B44:C44").Merge()

$ws.Range("B42:C42").UnMerge()
$ws.Range("B42:C42").Merge()

$ws.Range("A9:D9").UnMerge()
$ws.Range("A9:D9").Merge()
